$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.00841100000001
$ws.Range("H2").Value = 75.02523300000001
$ws.Range("I2").Value = 0.4156829172908309
$ws.Range("J2").Value = 0.415682917290831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 55.908252
$ws.Range("N2").Value = 167.724756
$ws.Range("O2").Value = 0.6412441619121594
$ws.Range("P2").Value = 0.6412441619121594
$ws.Range("Q2").Value = 1398.176544307572
$ws.Range("R2").Value = 12583.58889876815
$ws.Range("S2").Value = 0.2665542439193603
$ws.Range("T2").Value = 0.2665542439193604

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.00841100000001
$ws.Range("H3").Value = 75.02523300000001
$ws.Range("I3").Value = 0.4156829172908309
$ws.Range("J3").Value = 0.415682917290831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.370676
$ws.Range("N3").Value = 16.112028
$ws.Range("O3").Value = 0.06159939735768789
$ws.Range("P3").Value = 0.06159939735768789
$ws.Range("Q3").Value = 134.3120727558361
$ws.Range("R3").Value = 1208.808654802524
$ws.Range("S3").Value = 0.0256058171970008
$ws.Range("T3").Value = 0.02560581719700081

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.00841100000001
$ws.Range("H4").Value = 75.02523300000001
$ws.Range("I4").Value = 0.4156829172908309
$ws.Range("J4").Value = 0.415682917290831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.90822366666667
$ws.Range("N4").Value = 77.724671
$ws.Range("O4").Value = 0.2971564407301527
$ws.Range("P4").Value = 0.2971564407301527
$ws.Range("Q4").Value = 647.9235057359272
$ws.Range("R4").Value = 5831.311551623344
$ws.Range("S4").Value = 0.1235228561744698
$ws.Range("T4").Value = 0.1235228561744698

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.496019
$ws.Range("H5").Value = 61.488057
$ws.Range("I5").Value = 0.340679180727168
$ws.Range("J5").Value = 0.3406791807271681
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 55.908252
$ws.Range("N5").Value = 167.724756
$ws.Range("O5").Value = 0.6412441619121594
$ws.Range("P5").Value = 0.6412441619121594
$ws.Range("Q5").Value = 1145.896595248788
$ws.Range("R5").Value = 10313.06935723909
$ws.Range("S5").Value = 0.2184585357263139
$ws.Range("T5").Value = 0.218458535726314

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.496019
$ws.Range("H6").Value = 61.488057
$ws.Range("I6").Value = 0.340679180727168
$ws.Range("J6").Value = 0.3406791807271681
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.370676
$ws.Range("N6").Value = 16.112028
$ws.Range("O6").Value = 0.06159939735768789
$ws.Range("P6").Value = 0.06159939735768789
$ws.Range("Q6").Value = 110.077477338844
$ws.Range("R6").Value = 990.697296049596
$ws.Range("S6").Value = 0.02098563222510439
$ws.Range("T6").Value = 0.02098563222510439

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.496019
$ws.Range("H7").Value = 61.488057
$ws.Range("I7").Value = 0.340679180727168
$ws.Range("J7").Value = 0.3406791807271681
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 25.90822366666667
$ws.Range("N7").Value = 77.724671
$ws.Range("O7").Value = 0.2971564407301527
$ws.Range("P7").Value = 0.2971564407301527
$ws.Range("Q7").Value = 531.0154445282498
$ws.Range("R7").Value = 4779.139000754247
$ws.Range("S7").Value = 0.1012350127757497
$ws.Range("T7").Value = 0.1012350127757497

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.65779933333333
$ws.Range("H8").Value = 43.973398
$ws.Range("I8").Value = 0.243637901982001
$ws.Range("J8").Value = 0.243637901982001
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 55.908252
$ws.Range("N8").Value = 167.724756
$ws.Range("O8").Value = 0.6412441619121594
$ws.Range("P8").Value = 0.6412441619121594
$ws.Range("Q8").Value = 819.4919388934321
$ws.Range("R8").Value = 7375.427450040889
$ws.Range("S8").Value = 0.1562313822664851
$ws.Range("T8").Value = 0.1562313822664851

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.65779933333333
$ws.Range("H9").Value = 43.973398
$ws.Range("I9").Value = 0.243637901982001
$ws.Range("J9").Value = 0.243637901982001
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.370676
$ws.Range("N9").Value = 16.112028
$ws.Range("O9").Value = 0.06159939735768789
$ws.Range("P9").Value = 0.06159939735768789
$ws.Range("Q9").Value = 78.72229109234935
$ws.Range("R9").Value = 708.5006198311441
$ws.Range("S9").Value = 0.01500794793558269
$ws.Range("T9").Value = 0.0150079479355827

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.65779933333333
$ws.Range("H10").Value = 43.973398
$ws.Range("I10").Value = 0.243637901982001
$ws.Range("J10").Value = 0.243637901982001
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.90822366666667
$ws.Range("N10").Value = 77.724671
$ws.Range("O10").Value = 0.2971564407301527
$ws.Range("P10").Value = 0.2971564407301527
$ws.Range("Q10").Value = 379.7575435891176
$ws.Range("R10").Value = 3417.817892302058
$ws.Range("S10").Value = 0.07239857177993322
$ws.Range("T10").Value = 0.07239857177993324
